$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("H40").Value = 5311.125
$ws.Range("J40").Value = 5141.2856
$ws.Range("L40").Value = 5141.2856
$ws.Range("N40").Value = -5491.2856
$ws.Range("H92").Value = 883.76
$ws.Range("J92").Value = 636.7143
$ws.Range("L92").Value = 636.7143
$ws.Range("N92").Value = -3132.7143
$ws.Range("H94").Value = 823.1429000000001
$ws.Range("I94").Value = 823.1429000000001
$ws.Range("K94").Value = 823.1429000000001
$ws.Range("M94").Value = -372.1429000000001
$ws.Range("H98").Value = 2548.7932
$ws.Range("I98").Value = 2389.8215
$ws.Range("J98").Value = 7000
$ws.Range("K98").Value = 2389.8215
$ws.Range("L98").Value = 7000
$ws.Range("M98").Value = -891.8215
$ws.Range("N98").Value = -9996
$ws.Range("H122").Value = 2548.7932
$ws.Range("I122").Value = 2389.8215
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 7169.4645
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -4719.4645
$ws.Range("N122").Value = -25900
$ws.Range("H131").Value = 2012
$ws.Range("I131").Value = 2020.3334
$ws.Range("K131").Value = 6061.0002
$ws.Range("M131").Value = -1021.0002
$ws.Range("H132").Value = 37365.57
$ws.Range("I132").Value = 37365.57
$ws.Range("K132").Value = 112096.71
$ws.Range("M132").Value = -109566.71
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M19,M137,N137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1763377.2
$ws.Range("I32").Value = 1790390.8
$ws.Range("K32").Value = 1790390.8
$ws.Range("M32").Value = -1790103.8
$ws.Range("H43").Value = 26459.777
$ws.Range("J43").Value = 27306
$ws.Range("L43").Value = 27306
$ws.Range("N43").Value = -27932
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M74,N74,M77,N77").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 294568
$ws.Range("I99").Value = 11495.5
$ws.Range("K99").Value = 11495.5
$ws.Range("M99").Value = -9997.5
$ws.Range("H105").Value = 89665.95
$ws.Range("I105").Value = 2828.1765
$ws.Range("J105").Value = 335706.34
$ws.Range("K105").Value = 2828.1765
$ws.Range("L105").Value = 335706.34
$ws.Range("M105").Value = -1081.1765
$ws.Range("N105").Value = -339200.34

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 165
$ws.Range("I22").Value = 121.888885
$ws.Range("J22").Value = 208.11111
$ws.Range("K22").Value = 121.888885
$ws.Range("L22").Value = 208.11111
$ws.Range("M22").Value = 228.111115
$ws.Range("N22").Value = -908.1111100000001
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H99").Value = 3161.111
$ws.Range("I99").Value = 2608
$ws.Range("J99").Value = 3437.6667
$ws.Range("K99").Value = 2608
$ws.Range("L99").Value = 3437.6667
$ws.Range("M99").Value = -1110
$ws.Range("N99").Value = -6433.6667
$ws.Range("H100").Value = 780
$ws.Range("J100").Value = 780
$ws.Range("L100").Value = 780
$ws.Range("N100").Value = -2944
$ws.Range("H105").Value = 2472.1052
$ws.Range("I105").Value = 1852.7142
$ws.Range("J105").Value = 4206.4
$ws.Range("K105").Value = 1852.7142
$ws.Range("L105").Value = 4206.4
$ws.Range("M105").Value = -105.7141999999999
$ws.Range("N105").Value = -7700.4
$ws.Range("H122").Value = 2567.25
$ws.Range("I122").Value = 2567.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7701.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5251.75
$ws.Range("H126").Value = 3161.111
$ws.Range("I126").Value = 2608
$ws.Range("J126").Value = 3437.6667
$ws.Range("K126").Value = 7824
$ws.Range("L126").Value = 10313.0001
$ws.Range("M126").Value = -5354
$ws.Range("N126").Value = -15253.0001
$ws.Range("H132").Value = 2214.8215
$ws.Range("I132").Value = 2356.04
$ws.Range("J132").Value = 1038
$ws.Range("K132").Value = 7068.12
$ws.Range("L132").Value = 3114
$ws.Range("M132").Value = -4538.12
$ws.Range("N132").Value = -8174
$ws.Range("H134").Value = 4349530.5
$ws.Range("I134").Value = 1805.8422
$ws.Range("K134").Value = 5417.5266
$ws.Range("M134").Value = -2882.5266
$ws.Range("M31,N31,M34,N34,N122").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 5537.5
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("H134").Value = 660.125
$ws.Range("I134").Value = 660.125
$ws.Range("K134").Value = 1980.375
$ws.Range("M134").Value = 3089.625
$ws.Range("M47").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 6469.75
$ws.Range("I35").Value = 4999
$ws.Range("J35").Value = 6679.857
$ws.Range("K35").Value = 4999
$ws.Range("L35").Value = 6679.857
$ws.Range("M35").Value = -4701
$ws.Range("N35").Value = -7275.857
$ws.Range("H102").Value = 1131.25
$ws.Range("I102").Value = 1131.25
$ws.Range("K102").Value = 1131.25
$ws.Range("M102").Value = 490.75
$ws.Range("H117").Value = 26999
$ws.Range("J117").Value = 26999
$ws.Range("L117").Value = 26999
$ws.Range("N117").Value = -33883
$ws.Range("H126").Value = 6127.72
$ws.Range("I126").Value = 2291.5386
$ws.Range("K126").Value = 6874.6158
$ws.Range("M126").Value = -4404.6158
$ws.Range("H127").Value = 326
$ws.Range("J127").Value = 326
$ws.Range("L127").Value = 326
$ws.Range("N127").Value = -10246
$ws.Range("H132").Value = 1728.0385
$ws.Range("I132").Value = 1714.5652
$ws.Range("J132").Value = 1831.3334
$ws.Range("K132").Value = 5143.6956
$ws.Range("L132").Value = 5494.0002
$ws.Range("M132").Value = -2613.6956
$ws.Range("N132").Value = -10554.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3875
$ws.Range("I7").Value = 3458.5
$ws.Range("J7").Value = 4499.75
$ws.Range("K7").Value = 3458.5
$ws.Range("L7").Value = 4499.75
$ws.Range("M7").Value = -3346.5
$ws.Range("N7").Value = -4723.75
$ws.Range("H40").Value = 5198.7334
$ws.Range("I40").Value = 4900.4375
$ws.Range("K40").Value = 4900.4375
$ws.Range("M40").Value = -4764.4375
$ws.Range("H122").Value = 5630.6
$ws.Range("I122").Value = 2601.75
$ws.Range("J122").Value = 7649.8335
$ws.Range("K122").Value = 7805.25
$ws.Range("L122").Value = 22949.5005
$ws.Range("M122").Value = -5355.25
$ws.Range("N122").Value = -27849.5005
$ws.Range("H126").Value = 3875
$ws.Range("I126").Value = 3458.5
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 10375.5
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -7905.5
$ws.Range("N126").Value = -18439.25
$ws.Range("H136").Value = 29417304
$ws.Range("I136").Value = 5110.8965
$ws.Range("K136").Value = 15332.6895
$ws.Range("M136").Value = -12782.6895

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 29121.25
$ws.Range("J34").Value = 28828.334
$ws.Range("L34").Value = 28828.334
$ws.Range("N34").Value = -29234.334
$ws.Range("H62").Value = 5273
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 8003
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 8003
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -9251
$ws.Range("H65").Value = 5273
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 8003
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 40015
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -46255
$ws.Range("H96").Value = 2929.2778
$ws.Range("I96").Value = 2899.2
$ws.Range("J96").Value = 2940.8462
$ws.Range("K96").Value = 2899.2
$ws.Range("L96").Value = 2940.8462
$ws.Range("M96").Value = -1526.2
$ws.Range("N96").Value = -5686.8462
$ws.Range("H113").Value = 1129.7407
$ws.Range("I113").Value = 886.2727
$ws.Range("K113").Value = 2658.8181
$ws.Range("M113").Value = -488.8181
$ws.Range("H122").Value = 4714.6113
$ws.Range("I122").Value = 3774.8572
$ws.Range("J122").Value = 8003.75
$ws.Range("K122").Value = 11324.5716
$ws.Range("L122").Value = 24011.25
$ws.Range("M122").Value = -8874.571599999999
$ws.Range("N122").Value = -28911.25
$ws.Range("H126").Value = 1951.4333
$ws.Range("I126").Value = 1981.037
$ws.Range("K126").Value = 5943.111
$ws.Range("M126").Value = -3473.111
$ws.Range("H136").Value = 890.3333
$ws.Range("I136").Value = 855.19354
$ws.Range("K136").Value = 2565.58062
$ws.Range("M136").Value = -15.58061999999973
